$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 82: date serial 45884 (2025-08-15) and value 105.26
$ws.Range("A82").Value = 45884
$ws.Range("B82").Value = 105.26

# Copy the formatting (number format, font, border, alignment) from the
# row above (A81) so the new date cell matches the existing style.
$ws.Range("A81").Copy() | Out-Null
$ws.Range("A82").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
